$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprints")

# Row 4 - Sprint 1
$ws.Range("E4").Value = "6"
$ws.Range("F4").Value = "Choosing mini project theme. Sprint Planning document"
$ws.Range("F4").WrapText = $true
$ws.Range("G4").Value = 45796
$ws.Range("H4").Value = 45798
$ws.Range("I4").Value = "Closed"

# Row 5 - Sprint 2
$ws.Range("E5").Value = "6"
$ws.Range("F5").Value = "Version 1 of SRS. UML diagram, Functional REQs. Discussion on integration planning. Testing the development board. Assignment of workload"
$ws.Range("F5").WrapText = $true
$ws.Range("G5").Value = 45798
$ws.Range("H5").Value = 45801
$ws.Range("I5").Value = "Closed"

# Row 6 - Sprint 3
$ws.Range("E6").Value = "7"
$ws.Range("F6").Value = "SRS version 2. Website system architecture was defined. Further testing with the devlopment board. Sample codes were pushed into our respective branches."
$ws.Range("F6").WrapText = $true
$ws.Range("G6").Value = 45801
$ws.Range("H6").Value = 45809
$ws.Range("I6").Value = "Closed"

# Row 8 - Sprint 4
$ws.Range("G8").Value = 45857
$ws.Range("H8").Value = 45869
$ws.Range("I8").Value = "Active"

# Row 9 - Sprint 5 (keep formulas, values recompute naturally)
$ws.Range("G8").Value = 45857

# Row 10 - Sprint 6 (keep formulas)
